$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of numbers (23..51) across columns A:AC, continuing the sequence
# started in rows 2 and 5 (0-based index counters).
$vals = @(23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)

$newRow = 6
$rng = $ws.Range("A$newRow`:AC$newRow")

# Apply the new "40% lighter, Accent2-themed" fill to the whole new row in a
# single call so every cell in the row shares one style record.
$rng.Interior.Color = 8630772

for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($newRow, $i + 1).Value = $vals[$i]
}

# Move/collapse the selection onto S13, matching the post-edit cursor
# position saved in the workbook.
$ws.Range("S13").Select() | Out-Null
